$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.003.00"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.090.31"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.87"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0842"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "2.401.39"
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("E15").Value = "  +6.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.773"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").Value = "2.085.50"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").Value = "37.945.13"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +9.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.32%  "
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("E38").Value = "  +8.02%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("D41").Value = "1.544.49"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0905"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").Value = "2.288.84"
$ws.Range("E51").Value = "  +2.89%  "
